$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "updated at" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 16:35"

# 2. Row 4 - Estados Unidos
$ws.Range("B4").Value = 1648283
$ws.Range("C4").Value = 3189
$ws.Range("E4").Value = 1147239
$ws.Range("G4").Value = 85
$ws.Range("H4").Value = 97732

# 3. Row 11 - Alemania
$ws.Range("B11").Value = 179768
$ws.Range("C11").Value = 55
$ws.Range("E11").Value = 11514
$ws.Range("G11").Value = 2
$ws.Range("H11").Value = 8354

# 4. Rows 140/141 - Togo and Cabo Verde swap order (Cabo Verde moves above Togo)
# Row 140 becomes Cabo Verde with updated figures
$ws.Range("A140").Value = "Cabo Verde"
$ws.Range("B140").Value = 371
$ws.Range("C140").Value = 9
$ws.Range("D140").Value = 142
$ws.Range("E140").Value = 226
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 3

# Row 141 becomes Togo (figures unchanged from before, just relocated)
$ws.Range("A141").Value = "Togo"
$ws.Range("B141").Value = 363
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 121
$ws.Range("E141").Value = 230
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 12

# 5. Rows 149/150 - Santo Tome y Principe and Liberia swap order (Liberia moves above)
# Row 149 becomes Liberia with updated figures
$ws.Range("A149").Value = "Liberia"
$ws.Range("B149").Value = 255
$ws.Range("C149").Value = 6
$ws.Range("D149").Value = 136
$ws.Range("E149").Value = 93
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 2
$ws.Range("H149").Value = 26

# Row 150 becomes Santo Tome y Principe (figures unchanged from before, just relocated)
$ws.Range("A150").Value = "Santo Tome y Principe"
$ws.Range("B150").Value = 251
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 4
$ws.Range("E150").Value = 239
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 8
